# Update countries & provincias Spain
# - Re-sort two pairs of same-rank entries (Guinea / Consejo Danes para los
#   Refugiados, and Nueva Caledonia / Santa Lucia) by swapping their
#   displayed country-name text while leaving each row's own stats in place
#   (the underlying numeric stats for those rows are updated separately,
#   matching the refreshed data feed).
# - Refresh the "Datos actualizados" timestamp.
# - Refresh case/recovered/death counters for the rows whose source data
#   changed between the 13:32 and 14:49 snapshots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name reordering -------------------------------------------
# Row 99 held "Guinea" and row 100 held "Consejo Danes para los
# Refugiados"; the refreshed feed sorts them the other way around, so swap
# the two labels (use temporary placeholders so the two values actually
# exchange places instead of Excel treating the second write as a no-op).
$ws.Range("A99").Value = "__swap_placeholder_1__"
$ws.Range("A100").Value = "__swap_placeholder_2__"
$ws.Range("A99").Value = "Consejo Danes para los Refugiados"
$ws.Range("A100").Value = "Guinea"

# Row 207 held "Nueva Caledonia" and row 208 held "Santa Lucia"; swap them
# the same way.
$ws.Range("A207").Value = "__swap_placeholder_3__"
$ws.Range("A208").Value = "__swap_placeholder_4__"
$ws.Range("A207").Value = "Santa Lucia"
$ws.Range("A208").Value = "Nueva Caledonia"

# --- Timestamp -----------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 30 de Septiembre de 2020 a las 14:49"

# --- Refreshed counters ---------------------------------------------------
# Row 4 (Rusia)
$ws.Range("B4").Value = 7407201
$ws.Range("C4").Value = 1055
$ws.Range("D4").Value = 4649820
$ws.Range("E4").Value = 2546567
$ws.Range("G4").Value = 29
$ws.Range("H4").Value = 210814

# Row 5 (Colombia)
$ws.Range("B5").Value = 6233700
$ws.Range("C5").Value = 10181
$ws.Range("D5").Value = 5192917
$ws.Range("E5").Value = 943159
$ws.Range("G5").Value = 95
$ws.Range("H5").Value = 97624

# Row 19
$ws.Range("B19").Value = 362981
$ws.Range("C19").Value = 4691
$ws.Range("D19").Value = 292197
$ws.Range("E19").Value = 61603
$ws.Range("G19").Value = 59
$ws.Range("H19").Value = 9181

# Row 20
$ws.Range("B20").Value = 334605
$ws.Range("C20").Value = 418
$ws.Range("D20").Value = 319154
$ws.Range("E20").Value = 10683
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = 4768

# Row 25
$ws.Range("B25").Value = 290695
$ws.Range("C25").Value = 229
$ws.Range("E25").Value = 25139

# Row 35
$ws.Range("B35").Value = 120845
$ws.Range("C35").Value = 3294
$ws.Range("G35").Value = 13
$ws.Range("H35").Value = 6406

# Row 40
$ws.Range("B40").Value = 105182
$ws.Range("C40").Value = 614
$ws.Range("D40").Value = 96688
$ws.Range("E40").Value = 7884
$ws.Range("G40").Value = 3
$ws.Range("H40").Value = 610

# Row 44
$ws.Range("B44").Value = 92863
$ws.Range("G44").Value = 5
$ws.Range("H44").Value = 5893

# Row 49
$ws.Range("B49").Value = 78631
$ws.Range("C49").Value = 371
$ws.Range("D49").Value = 74525
$ws.Range("E49").Value = 3273
$ws.Range("G49").Value = 5
$ws.Range("H49").Value = 833

# Row 70
$ws.Range("B70").Value = 39899
$ws.Range("C70").Value = 358
$ws.Range("D70").Value = 31743
$ws.Range("E70").Value = 7845
$ws.Range("G70").Value = 5
$ws.Range("H70").Value = 311

# Row 75
$ws.Range("B75").Value = 34525
$ws.Range("C75").Value = 511
$ws.Range("D75").Value = 19361
$ws.Range("E75").Value = 14613
$ws.Range("G75").Value = 11
$ws.Range("H75").Value = 551

# Row 77
$ws.Range("B77").Value = 29077
$ws.Range("C77").Value = 96
$ws.Range("D77").Value = 23813
$ws.Range("E77").Value = 4421

# Row 78
$ws.Range("B78").Value = 27998
$ws.Range("C78").Value = 534
$ws.Range("D78").Value = 20754
$ws.Range("E78").Value = 6594

# Row 99 (now "Consejo Danes para los Refugiados")
$ws.Range("B99").Value = 10659
$ws.Range("C99").Value = 28
$ws.Range("D99").Value = 10139
$ws.Range("E99").Value = 248
$ws.Range("H99").Value = 272

# Row 100 (now "Guinea")
$ws.Range("B100").Value = 10634
$ws.Range("D100").Value = 9960
$ws.Range("E100").Value = 608
$ws.Range("H100").Value = 66

# Row 105
$ws.Range("D105").Value = 8100
$ws.Range("E105").Value = 1548
$ws.Range("H105").Value = 344

# Row 112
$ws.Range("B112").Value = 8129
$ws.Range("C112").Value = 112
$ws.Range("D112").Value = 4260
$ws.Range("E112").Value = 3794

# Row 150
$ws.Range("B150").Value = 2357
$ws.Range("C150").Value = 17
$ws.Range("D150").Value = 1973
$ws.Range("E150").Value = 343
$ws.Range("G150").Value = 1
$ws.Range("H150").Value = 41
